$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.947.83'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '2.978.12'
$ws.Range("E3").Value = '  -1.84%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '501.27'
$ws.Range("E5").Value = '  -3.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.91'
$ws.Range("E6").Value = '  -3.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.430'
$ws.Range("E8").Value = '  -2.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.30'
$ws.Range("E9").Value = '  -4.13%  '
$ws.Range("E10").Value = '  -2.17%  '
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").Value = '3.482.13'
$ws.Range("E12").Value = '  -1.92%  '
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.05'
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("E15").Value = '  -1.76%  '
$ws.Range("D16").Value = '57.027.73'
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.05'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '2.984.04'
$ws.Range("E18").Value = '  -2.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.61'
$ws.Range("E19").Value = '  -1.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.88'
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.46'
$ws.Range("E21").Value = '  -4.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.71'
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("E24").Value = '  -0.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.74'
$ws.Range("E25").Value = '  -1.38%  '
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("E27").Value = '  -5.41%  '
$ws.Range("D28").Value = '0.0₃0896'
$ws.Range("E28").Value = '  -4.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.54'
$ws.Range("E29").Value = '  -4.71%  '
$ws.Range("E30").Value = '  -1.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.77'
$ws.Range("E31").Value = '  -3.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.17'
$ws.Range("E32").Value = '  -4.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.18'
$ws.Range("E33").Value = '  -3.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '155.33'
$ws.Range("E34").Value = '  -1.69%  '
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("E37").Value = '  -4.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.02'
$ws.Range("E38").Value = '  -3.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0669'
$ws.Range("E39").Value = '  -2.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.71'
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = '3.010.01'
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.75'
$ws.Range("E43").Value = '  -0.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.639'
$ws.Range("E44").Value = '  -2.70%  '
$ws.Range("E45").Value = '  -3.94%  '
$ws.Range("D46").Value = '2.198.19'
$ws.Range("E46").Value = '  -5.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.946'
$ws.Range("E47").Value = '  -7.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.95'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E49").Value = '  -3.73%  '
$ws.Range("E50").Value = '  -2.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.82'
$ws.Range("E51").Value = '  -10.38%  '
